# Auto-generated Excel COM-interop script applying the Asura_Profits profit recalculation diff.
# For each affected row (identified by sheet name + row number), update columns H-N
# to reflect the new computed profit figures. Cells that are removed in the target
# state are cleared; cells that are newly introduced are written for the first time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 243.7619
$ws.Range("I33").Value = 237.73685
$ws.Range("J33").Value = 301
$ws.Range("K33").Value = 237.73685
$ws.Range("L33").Value = 301
$ws.Range("M33").Value = -8.736850000000004
$ws.Range("N33").Value = -759

$ws.Range("H62").Value = 2342.7778
$ws.Range("I62").Value = 2212.2144
$ws.Range("J62").Value = 2799.75
$ws.Range("K62").Value = 2212.2144
$ws.Range("L62").Value = 2799.75
$ws.Range("M62").Value = -1588.2144
$ws.Range("N62").Value = -4047.75

$ws.Range("H65").Value = 2342.7778
$ws.Range("I65").Value = 2212.2144
$ws.Range("J65").Value = 2799.75
$ws.Range("K65").Value = 11061.072
$ws.Range("L65").Value = 13998.75
$ws.Range("M65").Value = -7941.072
$ws.Range("N65").Value = -20238.75

$ws.Range("H70").Value = 11179056
$ws.Range("I70").Value = 33534232
$ws.Range("J70").Value = 1467.3
$ws.Range("K70").Value = 100602696
$ws.Range("L70").Value = 4401.9
$ws.Range("M70").Value = -100602426
$ws.Range("N70").Value = -4941.9

$ws.Range("H73").Value = 11179056
$ws.Range("I73").Value = 33534232
$ws.Range("J73").Value = 1467.3
$ws.Range("K73").Value = 100602696
$ws.Range("L73").Value = 4401.9
$ws.Range("M73").Value = -100601760
$ws.Range("N73").Value = -6273.9

$ws.Range("H106").Value = 777
$ws.Range("I106").Value = 471.25
$ws.Range("K106").Value = 471.25
$ws.Range("M106").Value = 159.75

$ws.Range("H112").Value = 1721.4828
$ws.Range("J112").Value = 1721.4828
$ws.Range("L112").Value = 5164.4484
$ws.Range("N112").Value = -7380.4484

$ws.Range("H132").Value = 2283.725
$ws.Range("I132").Value = 1854.5769
$ws.Range("J132").Value = 3080.7144
$ws.Range("K132").Value = 5563.7307
$ws.Range("L132").Value = 9242.143199999999
$ws.Range("M132").Value = -3033.7307
$ws.Range("N132").Value = -14302.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17895.299
$ws.Range("I32").Value = 18699.146
$ws.Range("J32").Value = 12402.333
$ws.Range("K32").Value = 18699.146
$ws.Range("L32").Value = 12402.333
$ws.Range("M32").Value = -18412.146
$ws.Range("N32").Value = -12976.333

$ws.Range("H132").Value = 1681.5333
$ws.Range("I132").Value = 1141.8096
$ws.Range("J132").Value = 2940.889
$ws.Range("K132").Value = 3425.4288
$ws.Range("L132").Value = 8822.667000000001
$ws.Range("M132").Value = -895.4288000000001
$ws.Range("N132").Value = -13882.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 98770
$ws.Range("J57").Value = 98770
$ws.Range("L57").Value = 98770
$ws.Range("N57").Value = -100210

$ws.Range("H122").Value = 70352.30499999999
$ws.Range("J122").Value = 70352.30499999999
$ws.Range("L122").Value = 70352.30499999999
$ws.Range("N122").Value = -80152.30499999999

$ws.Range("H134").Value = 24392250
$ws.Range("I134").Value = 27028800
$ws.Range("K134").Value = 81086400
$ws.Range("M134").Value = -81083865

$ws.Range("H136").Value = 98770
$ws.Range("J136").Value = 98770
$ws.Range("L136").Value = 98770
$ws.Range("N136").Value = -108970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25643984
$ws.Range("I31").Value = 37039144
$ws.Range("J31").Value = 4875
$ws.Range("K31").Value = 37039144
$ws.Range("L31").Value = 4875
$ws.Range("M31").Value = -37038849
$ws.Range("N31").Value = -5465

$ws.Range("H34").Value = 25643984
$ws.Range("I34").Value = 37039144
$ws.Range("J34").Value = 4875
$ws.Range("K34").Value = 37039144
$ws.Range("L34").Value = 4875
$ws.Range("M34").Value = -37038942
$ws.Range("N34").Value = -5279

$ws.Range("H134").Value = 1610.3077
$ws.Range("I134").Value = 1630.3636
$ws.Range("K134").Value = 4891.0908
$ws.Range("M134").Value = -2356.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 202
$ws.Range("I38").Value = 226.4
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 679.2
$ws.Range("L38").Value = 240
$ws.Range("M38").Value = -332.2
$ws.Range("N38").Value = -934

$ws.Range("H69").Value = 908.0909
$ws.Range("I69").Value = 430
$ws.Range("J69").Value = 1087.375
$ws.Range("K69").Value = 1290
$ws.Range("L69").Value = 3262.125
$ws.Range("M69").Value = -479
$ws.Range("N69").Value = -4884.125

$ws.Range("H72").Value = 908.0909
$ws.Range("I72").Value = 430
$ws.Range("J72").Value = 1087.375
$ws.Range("K72").Value = 3870
$ws.Range("L72").Value = 9786.375
$ws.Range("M72").Value = 186
$ws.Range("N72").Value = -17898.375

$ws.Range("H133").Value = 4270
$ws.Range("J133").Value = 6215
$ws.Range("L133").Value = 18645
$ws.Range("N133").Value = -28765

$ws.Range("H137").Value = 37039476
$ws.Range("I137").Value = 950
$ws.Range("J137").Value = 55558740
$ws.Range("K137").Value = 2850
$ws.Range("L137").Value = 166676220
$ws.Range("M137").Value = 2250
$ws.Range("N137").Value = -166686420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H102").Value = 2134.4443
$ws.Range("I102").Value = 2118.0908
$ws.Range("J102").Value = 2206.4
$ws.Range("K102").Value = 2118.0908
$ws.Range("L102").Value = 2206.4
$ws.Range("M102").Value = -496.0907999999999
$ws.Range("N102").Value = -5450.4

$ws.Range("H132").Value = 3320.9443
$ws.Range("I132").Value = 2414.7144
$ws.Range("K132").Value = 7244.1432
$ws.Range("M132").Value = -4714.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 50000
$ws.Range("I50").Value = 50000
$ws.Range("K50").Value = 50000
$ws.Range("M50").Value = -49363

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H132").Value = 6321.148
$ws.Range("I132").Value = 6417.7144
$ws.Range("J132").Value = 5983.1665
$ws.Range("K132").Value = 19253.1432
$ws.Range("L132").Value = 17949.4995
$ws.Range("M132").Value = -16723.1432
$ws.Range("N132").Value = -23009.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3850
$ws.Range("I62").Value = 3850
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3850
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3226
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3850
$ws.Range("I65").Value = 3850
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 19250
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -16130
$ws.Range("N65").ClearContents()

$ws.Range("H81").Value = 51737.45
$ws.Range("I81").Value = 72790.64
$ws.Range("J81").Value = 2613.3333
$ws.Range("K81").Value = 145581.28
$ws.Range("L81").Value = 5226.6666
$ws.Range("M81").Value = -144520.28
$ws.Range("N81").Value = -7348.6666

$ws.Range("H84").Value = 51737.45
$ws.Range("I84").Value = 72790.64
$ws.Range("J84").Value = 2613.3333
$ws.Range("K84").Value = 727906.4
$ws.Range("L84").Value = 26133.333
$ws.Range("M84").Value = -722602.4
$ws.Range("N84").Value = -36741.333

$ws.Range("H94").Value = 42115
$ws.Range("I94").Value = 39900
$ws.Range("J94").Value = 44330
$ws.Range("K94").Value = 39900
$ws.Range("L94").Value = 44330
$ws.Range("M94").Value = -38999
$ws.Range("N94").Value = -46132

$ws.Range("H96").Value = 950
$ws.Range("I96").Value = 633.3333
$ws.Range("J96").Value = 1900
$ws.Range("K96").Value = 633.3333
$ws.Range("L96").Value = 1900
$ws.Range("M96").Value = 739.6667
$ws.Range("N96").Value = -4646

$ws.Range("H132").Value = 1459.4117
$ws.Range("I132").Value = 1285.0312
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 3855.0936
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -1325.0936
$ws.Range("N132").Value = -17808.5
